# Nuuchahnulth data progress: update data progress
# Applies the "Christmas" week rows (21-26) of progress data on the Data
# sheet, then leaves the Graphs tab selected (matching the author's saved
# view state) with the Data sheet scrolled/selected at O27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------
# Fill in the page-count tracking columns (J:N) and the derived formulas
# (O = percent complete, P = percent remaining) for rows 21-26, which
# previously only had columns B-I populated.
# ---------------------------------------------------------------------

$rows = @(
  @{ Row = 21; J = 652; K = 189; L = 183; M = 183; N = 183 },
  @{ Row = 22; J = 652; K = 189; L = 183; M = 183; N = 183 },
  @{ Row = 23; J = 652; K = 189; L = 183; M = 183; N = 183 },
  @{ Row = 24; J = 652; K = 189; L = 183; M = 183; N = 183 },
  @{ Row = 25; J = 652; K = 247; L = 183; M = 183; N = 183 },
  @{ Row = 26; J = 652; K = 314; L = 183; M = 183; N = 183 }
)

foreach ($r in $rows) {
  $n = $r.Row
  $ws.Range("J$n").Value = $r.J
  $ws.Range("K$n").Value = $r.K
  $ws.Range("L$n").Value = $r.L
  $ws.Range("M$n").Value = $r.M
  $ws.Range("N$n").Value = $r.N
  $ws.Range("O$n").Formula = "=SUM((`$J$n-`$D`$1), (`$K$n-`$D`$1), (`$L$n-`$D`$1), (`$M$n-`$D`$1), (`$N$n-`$D`$1))/(`$D`$3*5)"
  $ws.Range("P$n").Formula = "=1-O$n"
}

# ---------------------------------------------------------------------
# View-state: Data sheet ends up scrolled to/selected at O27, and the
# Graphs sheet (not Data) is the active tab when the workbook is saved.
# ---------------------------------------------------------------------

$ws.Range("O27").Select()

$graphs = $wb.Worksheets.Item("Graphs")
$graphs.Activate()
